$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.4532998755493539
$ws.Range("C2").Value = 0.6247489020384521
$ws.Range("D2").Value = 0.7098241526097382
$ws.Range("E2").Value = 0.7334881754064316
$ws.Range("B3").Value = 0.4728680731960298
$ws.Range("C3").Value = 0.6374184673032106
$ws.Range("D3").Value = 0.7155518372566629
$ws.Range("E3").Value = 0.7369496738843104
$ws.Range("B4").Value = 0.4930149867371041
$ws.Range("C4").Value = 0.6626129769022936
$ws.Range("D4").Value = 0.7396515998568782
$ws.Range("E4").Value = 0.7602560484723868
$ws.Range("B5").Value = 0.5775338682106308
$ws.Range("C5").Value = 0.7200122765090208
$ws.Range("D5").Value = 0.7668410534285477
$ws.Range("E5").Value = 0.7861785817351283
$ws.Range("B6").Value = 0.5646114238160951
$ws.Range("C6").Value = 0.7144241959553093
$ws.Range("D6").Value = 0.7644877422023012
$ws.Range("E6").Value = 0.7846734571989462
$ws.Range("B7").Value = 0.6222607695612322
$ws.Range("C7").Value = 0.7489480418434787
$ws.Range("D7").Value = 0.7856843859769028
$ws.Range("E7").Value = 0.7892787435559668
$ws.Range("B8").Value = 0.2492231869568524
$ws.Range("C8").Value = 0.4597712142751087
$ws.Range("D8").Value = 0.5948629548070007
$ws.Range("E8").Value = 0.709022198366438
$ws.Range("B9").Value = 0.5286358997867628
$ws.Range("C9").Value = 0.6780164160477735
$ws.Range("D9").Value = 0.7453565595481404
$ws.Range("E9").Value = 0.7651098808327893
$ws.Range("B10").Value = 0.5960745203152207
$ws.Range("C10").Value = 0.7369860057981762
$ws.Range("D10").Value = 0.7814954067393087
$ws.Range("E10").Value = 0.7863713395062196
$ws.Range("B11").Value = 0.5990363700050438
$ws.Range("C11").Value = 0.7389962251113296
$ws.Range("D11").Value = 0.7828152745360828
$ws.Range("E11").Value = 0.78756876434718
$ws.Range("B12").Value = 0.6112096612008782
$ws.Range("C12").Value = 0.7464344466847087
$ws.Range("D12").Value = 0.7856065710047644
$ws.Range("E12").Value = 0.7885821271855722
$ws.Range("B13").Value = 0.5945170037428875
$ws.Range("C13").Value = 0.7356696593475156
$ws.Range("D13").Value = 0.7805575362458407
$ws.Range("E13").Value = 0.7855739167611877